# Generate Report for Handoff
# Adds a new localization-tracking entry for file
# "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" to the Overview / zh-cn / de-de
# sheets, placed immediately above the existing
# "e0dedbf3-747d-4a26-9822-c10367072f15.md" entry (which is pushed down one row).

$wb = $excel.ActiveWorkbook

function Set-CellText($ws, $row, $col, $text) {
    $ws.Cells.Item($row, $col).Value2 = $text
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1): columns A=File Name, B=zh-cn, C=de-de,
# D=Latest Handoff Date. Existing row 8 holds the "e0dedbf3..." entry and a
# hyperlink on A8. We first copy that row down to row 9 (value + style +
# formatting, via Range.Copy which preserves the hyperlink-blue style), then
# re-create its hyperlink on A9, and finally overwrite row 8 with the new
# "3f94c6a7..." entry and a fresh hyperlink on A8.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Push the existing "e0dedbf3" row from row 8 down to row 9 (copies value + style).
$wsOverview.Cells.Item(8, 1).Copy($wsOverview.Cells.Item(9, 1))
$wsOverview.Cells.Item(8, 2).Copy($wsOverview.Cells.Item(9, 2))
$wsOverview.Cells.Item(8, 3).Copy($wsOverview.Cells.Item(9, 3))
$wsOverview.Cells.Item(8, 4).Copy($wsOverview.Cells.Item(9, 4))

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d5f2ff3a31f94a94258aee4ec97e4fecf128319e/e2e/e0dedbf3-747d-4a26-9822-c10367072f15.md",
    "",
    "",
    "e0dedbf3-747d-4a26-9822-c10367072f15.md"
) | Out-Null

# Overwrite row 8 with the new "3f94c6a7" entry.
Set-CellText $wsOverview 8 1 "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md"
Set-CellText $wsOverview 8 2 "Ready for handoff"
Set-CellText $wsOverview 8 3 "Ready for handoff"
Set-CellText $wsOverview 8 4 "2016-28-11 16:28:14"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A8"),
    "https://github.com/OpenLocalizationTest/oltest/blob/2f6a6e5b1c0b99d8f1e3a4d7c8b2a5e6f7081930/e2e/3f94c6a7-95a4-4ba8-b751-690d346ef68e.md",
    "",
    "",
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Helper to populate a language sheet (zh-cn / de-de) row with the
# "Source File Name / File Extension / Status / Latest Handoff File /
#  Latest Handoff Datetime / ... / Latest Handback DateTime / Handoff Reason"
# columns used by both detail sheets (A..E, H, I).
# ---------------------------------------------------------------------------
function Set-DetailRow($ws, $row, $sourceMd, $xlfName, $handoffDate) {
    Set-CellText $ws $row 1 $sourceMd          # A: Source File Name
    Set-CellText $ws $row 2 ".md"              # B: File Extension
    Set-CellText $ws $row 3 "Ready for handoff" # C: Status
    Set-CellText $ws $row 4 $xlfName           # D: Latest Handoff File
    Set-CellText $ws $row 5 $handoffDate       # E: Latest Handoff Datetime
    Set-CellText $ws $row 8 "0001-01-01 00:00:00" # H: Latest Handback DateTime
    Set-CellText $ws $row 9 "Include"          # I: Handoff Reason
}

function Add-DetailHyperlinks($ws, $row, $mdUrl, $mdDisplay, $xlfUrl, $xlfDisplay) {
    $ws.Hyperlinks.Add($ws.Range("A" + $row), $mdUrl, "", "", $mdDisplay) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B" + $row), $mdUrl, "", "", ".md") | Out-Null
    $ws.Hyperlinks.Add($ws.Range("D" + $row), $xlfUrl, "", "", $xlfDisplay) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Push existing row 8 ("e0dedbf3...") down to row 9.
$wsZhCn.Cells.Item(8, 1).Copy($wsZhCn.Cells.Item(9, 1))
$wsZhCn.Cells.Item(8, 2).Copy($wsZhCn.Cells.Item(9, 2))
$wsZhCn.Cells.Item(8, 3).Copy($wsZhCn.Cells.Item(9, 3))
$wsZhCn.Cells.Item(8, 4).Copy($wsZhCn.Cells.Item(9, 4))
$wsZhCn.Cells.Item(8, 5).Copy($wsZhCn.Cells.Item(9, 5))
$wsZhCn.Cells.Item(8, 8).Copy($wsZhCn.Cells.Item(9, 8))
$wsZhCn.Cells.Item(8, 9).Copy($wsZhCn.Cells.Item(9, 9))

Add-DetailHyperlinks $wsZhCn 9 `
    "https://github.com/OpenLocalizationTest/oltest/blob/d5f2ff3a31f94a94258aee4ec97e4fecf128319e/e2e/e0dedbf3-747d-4a26-9822-c10367072f15.md" `
    "e0dedbf3-747d-4a26-9822-c10367072f15.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/17922edd66134fc971babdc4b4bacaff78b78964/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0dedbf3-747d-4a26-9822-c10367072f15.786df820bfbb0718682979bccfe4f64fc3e62dd7.zh-cn.xlf" `
    "e0dedbf3-747d-4a26-9822-c10367072f15.786df820bfbb0718682979bccfe4f64fc3e62dd7.zh-cn.xlf"

# Overwrite row 8 with the new "3f94c6a7" entry.
Set-DetailRow $wsZhCn 8 "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" `
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.0e466d4e78459a28a29d7a8ed1d5cdb37a0114e6.zh-cn.xlf" `
    "2016-03-11 16:28:11"

Add-DetailHyperlinks $wsZhCn 8 `
    "https://github.com/OpenLocalizationTest/oltest/blob/2f6a6e5b1c0b99d8f1e3a4d7c8b2a5e6f7081930/e2e/3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" `
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c4d6e8a2b19f0c3d5e6a7b8c9d0e1f203142536/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3f94c6a7-95a4-4ba8-b751-690d346ef68e.0e466d4e78459a28a29d7a8ed1d5cdb37a0114e6.zh-cn.xlf" `
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.0e466d4e78459a28a29d7a8ed1d5cdb37a0114e6.zh-cn.xlf"

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Push existing row 8 ("e0dedbf3...") down to row 9.
$wsDeDe.Cells.Item(8, 1).Copy($wsDeDe.Cells.Item(9, 1))
$wsDeDe.Cells.Item(8, 2).Copy($wsDeDe.Cells.Item(9, 2))
$wsDeDe.Cells.Item(8, 3).Copy($wsDeDe.Cells.Item(9, 3))
$wsDeDe.Cells.Item(8, 4).Copy($wsDeDe.Cells.Item(9, 4))
$wsDeDe.Cells.Item(8, 5).Copy($wsDeDe.Cells.Item(9, 5))
$wsDeDe.Cells.Item(8, 8).Copy($wsDeDe.Cells.Item(9, 8))
$wsDeDe.Cells.Item(8, 9).Copy($wsDeDe.Cells.Item(9, 9))

Add-DetailHyperlinks $wsDeDe 9 `
    "https://github.com/OpenLocalizationTest/oltest/blob/d5f2ff3a31f94a94258aee4ec97e4fecf128319e/e2e/e0dedbf3-747d-4a26-9822-c10367072f15.md" `
    "e0dedbf3-747d-4a26-9822-c10367072f15.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5644dac920b3d4470e6c2b4b188e7324fc4c4b9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0dedbf3-747d-4a26-9822-c10367072f15.786df820bfbb0718682979bccfe4f64fc3e62dd7.de-de.xlf" `
    "e0dedbf3-747d-4a26-9822-c10367072f15.786df820bfbb0718682979bccfe4f64fc3e62dd7.de-de.xlf"

# Overwrite row 8 with the new "3f94c6a7" entry.
Set-DetailRow $wsDeDe 8 "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" `
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.0e466d4e78459a28a29d7a8ed1d5cdb37a0114e6.de-de.xlf" `
    "2016-03-11 16:28:14"

Add-DetailHyperlinks $wsDeDe 8 `
    "https://github.com/OpenLocalizationTest/oltest/blob/2f6a6e5b1c0b99d8f1e3a4d7c8b2a5e6f7081930/e2e/3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" `
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9d2e4f6a8b0c1d3e5f7081930a2b4c6d8e0f1a2b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3f94c6a7-95a4-4ba8-b751-690d346ef68e.0e466d4e78459a28a29d7a8ed1d5cdb37a0114e6.de-de.xlf" `
    "3f94c6a7-95a4-4ba8-b751-690d346ef68e.0e466d4e78459a28a29d7a8ed1d5cdb37a0114e6.de-de.xlf"

Write-Output "Done."
